$d = $word.ActiveDocument

# 1. library(openxlsx) -> library(openxlsx); — add trailing semicolon
$d.Content.Find.Execute("(openxlsx)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(openxlsx);", 2)

# 2. Collapse the comment + line break + View(data) run sequence into a
#    single tightened comment, removing the extra "View(data)" line.
$d.Content.Find.Execute("# read the simple file^lView(data)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "#read the simple file", 2)
